$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "aaaa"
$ws.Range("B3").Value = "bbbb"
$ws.Range("B4").Value = "dddd"
$ws.Range("B5").Value = "github for win2"
$ws.Range("B6").Value = "wn3"
$ws.Range("B7").Value = "bbb"

$ws.Range("B8").Select()
